$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) to the table, mirroring the existing
# per-row formatting that's already used by column P.

function Set-Q {
    param(
        $Sheet,
        [int]$Row,
        $Val
    )
    $Sheet.Range("P$Row").Copy() | Out-Null
    $Sheet.Range("Q$Row").PasteSpecial(-4122) | Out-Null
    $Sheet.Range("Q$Row").Value = $Val
}

Set-Q $ws 4  2020
Set-Q $ws 5  0.1
Set-Q $ws 6  0.2
Set-Q $ws 7  "-"
Set-Q $ws 8  0.2
Set-Q $ws 9  "-"
Set-Q $ws 10 0.1
Set-Q $ws 11 "-"
Set-Q $ws 12 0.3
Set-Q $ws 13 "-"
Set-Q $ws 14 "-"

$excel.CutCopyMode = $false

$ws.Range("O17").Select() | Out-Null
